$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.175.59"
$ws.Range("E2").Value = "  -3.93%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.491.31"
$ws.Range("E3").Value = "  -5.34%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.82"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.54"
$ws.Range("E6").Value = "  -8.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.492.93"
$ws.Range("E7").Value = "  -5.20%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.481"
$ws.Range("E9").Value = "  -3.57%  "
$ws.Range("E10").Value = "  -5.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.91"
$ws.Range("E11").Value = "  -3.60%  "
$ws.Range("E12").Value = "  -5.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000214"
$ws.Range("E13").Value = "  -7.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.076.62"
$ws.Range("E14").Value = "  -5.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.41"
$ws.Range("E15").Value = "  -4.26%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.075.19"
$ws.Range("E16").Value = "  -3.94%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.474.32"
$ws.Range("E17").Value = "  -5.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.117"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.44"
$ws.Range("E19").Value = "  -1.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.05"
$ws.Range("E20").Value = "  -5.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "443.89"
$ws.Range("E21").Value = "  -5.90%  "
$ws.Range("E22").Value = "  -12.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.624"
$ws.Range("E23").Value = "  -5.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.13"
$ws.Range("E24").Value = "  -3.42%  "
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000126"
$ws.Range("E26").Value = "  -1.44%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.625.93"
$ws.Range("E27").Value = "  -5.34%  "
$ws.Range("E28").Value = "  -9.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.30"
$ws.Range("E29").Value = "  -6.75%  "
$ws.Range("E30").Value = "  -4.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.58"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("E33").Value = "  -1.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.63"
$ws.Range("E34").Value = "  -4.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.13"
$ws.Range("E35").Value = "  -6.38%  "
$ws.Range("E36").Value = "  -7.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.468.27"
$ws.Range("E37").Value = "  -5.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.98"
$ws.Range("E38").Value = "  -5.26%  "
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "169.41"
$ws.Range("E42").Value = "  -4.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0869"
$ws.Range("E43").Value = "  -3.62%  "
$ws.Range("E44").Value = "  -7.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.882"
$ws.Range("E45").Value = "  -5.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.64"
$ws.Range("E46").Value = "  -2.10%  "
$ws.Range("E47").Value = "  -0.74%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.25"
$ws.Range("E48").Value = "  -9.77%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.52"
$ws.Range("E49").Value = "  -11.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.56"
$ws.Range("E50").Value = "  -4.26%  "
$ws.Range("E51").Value = "  -4.44%  "
